$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 17709
$ws.Range("J87").Value = 17709
$ws.Range("L87").Value = 17709
$ws.Range("N87").Value = -20205
$ws.Range("H90").Value = 17709
$ws.Range("J90").Value = 17709
$ws.Range("L90").Value = 53127
$ws.Range("N90").Value = -65607
$ws.Range("H98").Value = 31103.424
$ws.Range("I98").Value = 1291.5
$ws.Range("J98").Value = 195069
$ws.Range("K98").Value = 1291.5
$ws.Range("L98").Value = 195069
$ws.Range("M98").Value = 206.5
$ws.Range("N98").Value = -198065
$ws.Range("H111").Value = 1075
$ws.Range("I111").Value = 1100
$ws.Range("J111").Value = 1050
$ws.Range("K111").Value = 3300
$ws.Range("L111").Value = 3150
$ws.Range("M111").Value = -233
$ws.Range("N111").Value = -9284
$ws.Range("H122").Value = 31103.424
$ws.Range("I122").Value = 1291.5
$ws.Range("J122").Value = 195069
$ws.Range("K122").Value = 3874.5
$ws.Range("L122").Value = 585207
$ws.Range("M122").Value = -1424.5
$ws.Range("N122").Value = -590107
$ws.Range("H126").Value = 46766
$ws.Range("J126").Value = 46766
$ws.Range("L126").Value = 46766
$ws.Range("N126").Value = -56646
$ws.Range("H130").Value = 49772
$ws.Range("J130").Value = 49772
$ws.Range("L130").Value = 49772
$ws.Range("N130").Value = -59812
$ws.Range("H131").Value = 1807.6875
$ws.Range("J131").Value = 2047
$ws.Range("L131").Value = 6141
$ws.Range("N131").Value = -16221
$ws.Range("H132").Value = 14624.116
$ws.Range("I132").Value = 2132.7288
$ws.Range("J132").Value = 88323.3
$ws.Range("K132").Value = 6398.1864
$ws.Range("L132").Value = 264969.9
$ws.Range("M132").Value = -3868.1864
$ws.Range("N132").Value = -270029.9
$ws.Range("H137").Value = 3810.1667
$ws.Range("I137").Value = 1092.2333
$ws.Range("J137").Value = 8340.056
$ws.Range("K137").Value = 3276.699900000001
$ws.Range("L137").Value = 25020.168
$ws.Range("M137").Value = -726.6999000000005
$ws.Range("N137").Value = -30120.168
$ws.Range("H138").Value = 1552.26
$ws.Range("I138").Value = 777.87177
$ws.Range("J138").Value = 2047.3606
$ws.Range("K138").Value = 2333.61531
$ws.Range("L138").Value = 6142.0818
$ws.Range("M138").Value = 2806.38469
$ws.Range("N138").Value = -16422.0818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 40369
$ws.Range("J80").Value = 40369
$ws.Range("L80").Value = 40369
$ws.Range("N80").Value = -42365
$ws.Range("H83").Value = 40369
$ws.Range("J83").Value = 40369
$ws.Range("L83").Value = 121107
$ws.Range("N83").Value = -131091
$ws.Range("H110").Value = 1599.8823
$ws.Range("I110").Value = 1683.7333
$ws.Range("J110").Value = 971
$ws.Range("K110").Value = 1683.7333
$ws.Range("L110").Value = 971
$ws.Range("M110").Value = 361.2666999999999
$ws.Range("N110").Value = -5061
$ws.Range("H123").Value = 40221.8
$ws.Range("J123").Value = 40221.8
$ws.Range("L123").Value = 40221.8
$ws.Range("N123").Value = -50021.8
$ws.Range("H132").Value = 8622087
$ws.Range("I132").Value = 12500899
$ws.Range("J132").Value = 2503.3333
$ws.Range("K132").Value = 37502697
$ws.Range("L132").Value = 7509.999899999999
$ws.Range("M132").Value = -37500167
$ws.Range("N132").Value = -12569.9999
$ws.Range("H135").Value = 35788.1
$ws.Range("J135").Value = 35788.1
$ws.Range("L135").Value = 35788.1
$ws.Range("N135").Value = -45928.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3163.5535
$ws.Range("I134").Value = 1402.625
$ws.Range("J134").Value = 3867.925
$ws.Range("K134").Value = 4207.875
$ws.Range("L134").Value = 11603.775
$ws.Range("M134").Value = -1672.875
$ws.Range("N134").Value = -16673.775
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 47157.5
$ws.Range("J81").Value = 47157.5
$ws.Range("L81").Value = 47157.5
$ws.Range("N81").Value = -49153.5
$ws.Range("H82").Value = 44173
$ws.Range("J82").Value = 44173
$ws.Range("L82").Value = 44173
$ws.Range("N82").Value = -44895
$ws.Range("H84").Value = 47157.5
$ws.Range("J84").Value = 47157.5
$ws.Range("L84").Value = 141472.5
$ws.Range("N84").Value = -151456.5
$ws.Range("H85").Value = 44173
$ws.Range("J85").Value = 44173
$ws.Range("L85").Value = 44173
$ws.Range("N85").Value = -46669
$ws.Range("H88").Value = 40748.75
$ws.Range("J88").Value = 40748.75
$ws.Range("L88").Value = 40748.75
$ws.Range("N88").Value = -41560.75
$ws.Range("H91").Value = 40748.75
$ws.Range("J91").Value = 40748.75
$ws.Range("L91").Value = 40748.75
$ws.Range("N91").Value = -43556.75
$ws.Range("H100").Value = 35913
$ws.Range("J100").Value = 35913
$ws.Range("L100").Value = 35913
$ws.Range("N100").Value = -38077
$ws.Range("H107").Value = 1007.3
$ws.Range("I107").Value = 996.625
$ws.Range("J107").Value = 1050
$ws.Range("K107").Value = 996.625
$ws.Range("L107").Value = 1050
$ws.Range("M107").Value = 923.375
$ws.Range("N107").Value = -4890
$ws.Range("H132").Value = 38538.633
$ws.Range("I132").Value = 1278.6333
$ws.Range("J132").Value = 178263.62
$ws.Range("K132").Value = 3835.8999
$ws.Range("L132").Value = 534790.86
$ws.Range("M132").Value = -1305.8999
$ws.Range("N132").Value = -539850.86
$ws.Range("H134").Value = 425870.78
$ws.Range("I134").Value = 1468.1111
$ws.Range("J134").Value = 2335682.8
$ws.Range("K134").Value = 4404.3333
$ws.Range("L134").Value = 7007048.399999999
$ws.Range("M134").Value = -1869.3333
$ws.Range("N134").Value = -7012118.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 10947.263
$ws.Range("I107").Value = 10384.6
$ws.Range("J107").Value = 11572.444
$ws.Range("K107").Value = 31153.8
$ws.Range("L107").Value = 34717.33199999999
$ws.Range("M107").Value = -29233.8
$ws.Range("N107").Value = -38557.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 47399.5
$ws.Range("J104").Value = 47399.5
$ws.Range("L104").Value = 47399.5
$ws.Range("N104").Value = -54387.5
$ws.Range("H107").Value = 2819.6
$ws.Range("I107").Value = 316.4
$ws.Range("J107").Value = 4071.2
$ws.Range("K107").Value = 316.4
$ws.Range("L107").Value = 4071.2
$ws.Range("M107").Value = 1603.6
$ws.Range("N107").Value = -7911.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2027.125
$ws.Range("I61").Value = 2101.36
$ws.Range("J61").Value = 1762
$ws.Range("K61").Value = 2101.36
$ws.Range("L61").Value = 1762
$ws.Range("M61").Value = -1899.36
$ws.Range("N61").Value = -2166
$ws.Range("H81").Value = 26500
$ws.Range("J81").Value = 26500
$ws.Range("L81").Value = 26500
$ws.Range("N81").Value = -28496
$ws.Range("H84").Value = 26500
$ws.Range("J84").Value = 26500
$ws.Range("L84").Value = 79500
$ws.Range("N84").Value = -89484
$ws.Range("H86").Value = 46191
$ws.Range("J86").Value = 46191
$ws.Range("L86").Value = 46191
$ws.Range("N86").Value = -48563
$ws.Range("H88").Value = 43178.332
$ws.Range("J88").Value = 43178.332
$ws.Range("L88").Value = 43178.332
$ws.Range("N88").Value = -44034.332
$ws.Range("H89").Value = 46191
$ws.Range("J89").Value = 46191
$ws.Range("L89").Value = 138573
$ws.Range("N89").Value = -150429
$ws.Range("H91").Value = 43178.332
$ws.Range("J91").Value = 43178.332
$ws.Range("L91").Value = 43178.332
$ws.Range("N91").Value = -46142.332
$ws.Range("H92").Value = 45381
$ws.Range("J92").Value = 45381
$ws.Range("L92").Value = 45381
$ws.Range("N92").Value = -50373
$ws.Range("H99").Value = 23610
$ws.Range("I99").Value = 12220
$ws.Range("K99").Value = 12220
$ws.Range("M99").Value = -9225
$ws.Range("H113").Value = 2027.125
$ws.Range("I113").Value = 2101.36
$ws.Range("J113").Value = 1762
$ws.Range("K113").Value = 2101.36
$ws.Range("L113").Value = 1762
$ws.Range("M113").Value = 68.63999999999987
$ws.Range("N113").Value = -6102
$ws.Range("H140").Value = 29285.6
$ws.Range("J140").Value = 29285.6
$ws.Range("L140").Value = 29285.6
$ws.Range("N140").Value = -39645.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 31333.334
$ws.Range("J92").Value = 31333.334
$ws.Range("L92").Value = 31333.334
$ws.Range("N92").Value = -36325.334
$ws.Range("H93").Value = 45000
$ws.Range("J93").Value = 45000
$ws.Range("L93").Value = 45000
$ws.Range("N93").Value = -49992
$ws.Range("H94").Value = 39913
$ws.Range("J94").Value = 39913
$ws.Range("L94").Value = 39913
$ws.Range("N94").Value = -41715
$ws.Range("H132").Value = 1631.8379
$ws.Range("I132").Value = 1383.3667
$ws.Range("J132").Value = 2696.7144
$ws.Range("K132").Value = 4150.1001
$ws.Range("L132").Value = 8090.1432
$ws.Range("M132").Value = -1620.1001
$ws.Range("N132").Value = -13150.1432
$ws.Range("H136").Value = 196807.02
$ws.Range("I136").Value = 222742.48
$ws.Range("J136").Value = 2291
$ws.Range("K136").Value = 668227.4400000001
$ws.Range("L136").Value = 6873
$ws.Range("M136").Value = -665677.4400000001
$ws.Range("N136").Value = -11973
